$wb = $excel.ActiveWorkbook

# --- Sheet "CPbE-FoCSbS" (CCS Percentages by Entity, fuel mix) ---
$ws3 = $wb.Worksheets.Item("CPbE-FoCSbS")
    $ws3.Range('B2').Formula = '=About!$I$24'
    $ws3.Range('C2').Formula = '=About!$I$24'
    $ws3.Range('D2').Formula = '=About!$I$24'
    $ws3.Range('E2').Formula = '=About!$I$24'
    $ws3.Range('F2').Formula = '=About!$I$24'
    $ws3.Range('G2').Formula = '=About!$I$24'
    $ws3.Range('H2').Formula = '=About!$I$24'
    $ws3.Range('I2').Formula = '=About!$I$24'
    $ws3.Range('J2').Formula = '=About!$I$24'
    $ws3.Range('K2').Formula = '=About!$I$24'
    $ws3.Range('L2').Formula = '=About!$I$24'
    $ws3.Range('M2').Formula = '=About!$I$24'
    $ws3.Range('N2').Formula = '=About!$I$24'
    $ws3.Range('O2').Formula = '=About!$I$24'
    $ws3.Range('P2').Formula = '=About!$I$24'
    $ws3.Range('Q2').Formula = '=About!$I$24'
    $ws3.Range('R2').Formula = '=About!$I$24'
    $ws3.Range('S2').Formula = '=About!$I$24'
    $ws3.Range('T2').Formula = '=About!$I$24'
    $ws3.Range('U2').Formula = '=About!$I$24'
    $ws3.Range('V2').Formula = '=About!$I$24'
    $ws3.Range('W2').Formula = '=About!$I$24'
    $ws3.Range('X2').Formula = '=About!$I$24'
    $ws3.Range('Y2').Formula = '=About!$I$24'
    $ws3.Range('Z2').Formula = '=About!$I$24'
    $ws3.Range('AA2').Formula = '=About!$I$24'
    $ws3.Range('AB2').Formula = '=About!$I$24'
    $ws3.Range('AC2').Formula = '=About!$I$24'
    $ws3.Range('AD2').Formula = '=About!$I$24'
    $ws3.Range('AE2').Formula = '=About!$I$24'
    $ws3.Range('AF2').Formula = '=About!$I$24'
    $ws3.Range('AG2').Formula = '=About!$I$24'
    $ws3.Range('AH2').Formula = '=About!$I$24'
    $ws3.Range('AI2').Formula = '=About!$I$24'
    $ws3.Range('AJ2').Formula = '=About!$I$24'
    $ws3.Range('AK2').Formula = '=About!$I$24'
    $ws3.Range('AL2').Formula = '=About!$I$24'
    $ws3.Range('AM2').Formula = '=About!$I$24'
    $ws3.Range('B3').Formula = '=About!$I$23'
# --- Sheet "CPbE-FoESCbES" (CCS Percentages by Entity, fuel+tech mix) ---
$ws4 = $wb.Worksheets.Item("CPbE-FoESCbES")
    $ws4.Range('B2').Formula = '=C2'
    $ws4.Range('C2').Formula = '=D2'
    $ws4.Range('D2').Formula = '=E2'
    $ws4.Range('E2').Formula = '=F2'
    $ws4.Range('F2').Formula = '=''Technologies to power CCS''!B19'
    $ws4.Range('G2').Formula = '=''Technologies to power CCS''!C19'
    $ws4.Range('H2').Formula = '=''Technologies to power CCS''!D19'
    $ws4.Range('I2').Formula = '=''Technologies to power CCS''!E19'
    $ws4.Range('J2').Formula = '=''Technologies to power CCS''!F19'
    $ws4.Range('K2').Formula = '=''Technologies to power CCS''!G19'
    $ws4.Range('L2').Formula = '=''Technologies to power CCS''!H19'
    $ws4.Range('M2').Formula = '=''Technologies to power CCS''!I19'
    $ws4.Range('N2').Formula = '=''Technologies to power CCS''!J19'
    $ws4.Range('O2').Formula = '=''Technologies to power CCS''!K19'
    $ws4.Range('P2').Formula = '=''Technologies to power CCS''!L19'
    $ws4.Range('Q2').Formula = '=''Technologies to power CCS''!M19'
    $ws4.Range('R2').Formula = '=''Technologies to power CCS''!N19'
    $ws4.Range('S2').Formula = '=''Technologies to power CCS''!O19'
    $ws4.Range('T2').Formula = '=S2'
    $ws4.Range('U2').Formula = '=T2'
    $ws4.Range('V2').Formula = '=U2'
    $ws4.Range('W2').Formula = '=V2'
    $ws4.Range('X2').Formula = '=W2'
    $ws4.Range('Y2').Formula = '=X2'
    $ws4.Range('Z2').Formula = '=Y2'
    $ws4.Range('AA2').Formula = '=Z2'
    $ws4.Range('AB2').Formula = '=AA2'
    $ws4.Range('AC2').Formula = '=AB2'
    $ws4.Range('AD2').Formula = '=AC2'
    $ws4.Range('AE2').Formula = '=AD2'
    $ws4.Range('AF2').Formula = '=AE2'
    $ws4.Range('AG2').Formula = '=AF2'
    $ws4.Range('AH2').Formula = '=AG2'
    $ws4.Range('AI2').Formula = '=AH2'
    $ws4.Range('AJ2').Formula = '=AI2'
    $ws4.Range('AK2').Formula = '=AJ2'
    $ws4.Range('AL2').Formula = '=AK2'
    $ws4.Range('AM2').Formula = '=AL2'
    $ws4.Range('B3').Formula = '=C3'
    $ws4.Range('C3').Formula = '=D3'
    $ws4.Range('D3').Formula = '=E3'
    $ws4.Range('E3').Formula = '=F3'
    $ws4.Range('F3').Formula = '=''Technologies to power CCS''!B20'
    $ws4.Range('G3').Formula = '=''Technologies to power CCS''!C20'
    $ws4.Range('H3').Formula = '=''Technologies to power CCS''!D20'
    $ws4.Range('I3').Formula = '=''Technologies to power CCS''!E20'
    $ws4.Range('J3').Formula = '=''Technologies to power CCS''!F20'
    $ws4.Range('K3').Formula = '=''Technologies to power CCS''!G20'
    $ws4.Range('L3').Formula = '=''Technologies to power CCS''!H20'
    $ws4.Range('M3').Formula = '=''Technologies to power CCS''!I20'
    $ws4.Range('N3').Formula = '=''Technologies to power CCS''!J20'
    $ws4.Range('O3').Formula = '=''Technologies to power CCS''!K20'
    $ws4.Range('P3').Formula = '=''Technologies to power CCS''!L20'
    $ws4.Range('Q3').Formula = '=''Technologies to power CCS''!M20'
    $ws4.Range('R3').Formula = '=''Technologies to power CCS''!N20'
    $ws4.Range('S3').Formula = '=''Technologies to power CCS''!O20'
    $ws4.Range('T3').Formula = '=S3'
    $ws4.Range('U3').Formula = '=T3'
    $ws4.Range('V3').Formula = '=U3'
    $ws4.Range('W3').Formula = '=V3'
    $ws4.Range('X3').Formula = '=W3'
    $ws4.Range('Y3').Formula = '=X3'
    $ws4.Range('Z3').Formula = '=Y3'
    $ws4.Range('AA3').Formula = '=Z3'
    $ws4.Range('AB3').Formula = '=AA3'
    $ws4.Range('AC3').Formula = '=AB3'
    $ws4.Range('AD3').Formula = '=AC3'
    $ws4.Range('AE3').Formula = '=AD3'
    $ws4.Range('AF3').Formula = '=AE3'
    $ws4.Range('AG3').Formula = '=AF3'
    $ws4.Range('AH3').Formula = '=AG3'
    $ws4.Range('AI3').Formula = '=AH3'
    $ws4.Range('AJ3').Formula = '=AI3'
    $ws4.Range('AK3').Formula = '=AJ3'
    $ws4.Range('AL3').Formula = '=AK3'
    $ws4.Range('AM3').Formula = '=AL3'
    $ws4.Range('B4').Formula = '=C4'
    $ws4.Range('C4').Formula = '=D4'
    $ws4.Range('D4').Formula = '=E4'
    $ws4.Range('E4').Formula = '=F4'
    $ws4.Range('F4').Formula = '=''Technologies to power CCS''!B21'
    $ws4.Range('G4').Formula = '=''Technologies to power CCS''!C21'
    $ws4.Range('H4').Formula = '=''Technologies to power CCS''!D21'
    $ws4.Range('I4').Formula = '=''Technologies to power CCS''!E21'
    $ws4.Range('J4').Formula = '=''Technologies to power CCS''!F21'
    $ws4.Range('K4').Formula = '=''Technologies to power CCS''!G21'
    $ws4.Range('L4').Formula = '=''Technologies to power CCS''!H21'
    $ws4.Range('M4').Formula = '=''Technologies to power CCS''!I21'
    $ws4.Range('N4').Formula = '=''Technologies to power CCS''!J21'
    $ws4.Range('O4').Formula = '=''Technologies to power CCS''!K21'
    $ws4.Range('P4').Formula = '=''Technologies to power CCS''!L21'
    $ws4.Range('Q4').Formula = '=''Technologies to power CCS''!M21'
    $ws4.Range('R4').Formula = '=''Technologies to power CCS''!N21'
    $ws4.Range('S4').Formula = '=''Technologies to power CCS''!O21'
    $ws4.Range('T4').Formula = '=S4'
    $ws4.Range('U4').Formula = '=T4'
    $ws4.Range('V4').Formula = '=U4'
    $ws4.Range('W4').Formula = '=V4'
    $ws4.Range('X4').Formula = '=W4'
    $ws4.Range('Y4').Formula = '=X4'
    $ws4.Range('Z4').Formula = '=Y4'
    $ws4.Range('AA4').Formula = '=Z4'
    $ws4.Range('AB4').Formula = '=AA4'
    $ws4.Range('AC4').Formula = '=AB4'
    $ws4.Range('AD4').Formula = '=AC4'
    $ws4.Range('AE4').Formula = '=AD4'
    $ws4.Range('AF4').Formula = '=AE4'
    $ws4.Range('AG4').Formula = '=AF4'
    $ws4.Range('AH4').Formula = '=AG4'
    $ws4.Range('AI4').Formula = '=AH4'
    $ws4.Range('AJ4').Formula = '=AI4'
    $ws4.Range('AK4').Formula = '=AJ4'
    $ws4.Range('AL4').Formula = '=AK4'
    $ws4.Range('AM4').Formula = '=AL4'
    $ws4.Range('B5').Formula = '=C5'
    $ws4.Range('C5').Formula = '=D5'
    $ws4.Range('D5').Formula = '=E5'
    $ws4.Range('E5').Formula = '=F5'
    $ws4.Range('F5').Formula = '=''Technologies to power CCS''!B22'
    $ws4.Range('G5').Formula = '=''Technologies to power CCS''!C22'
    $ws4.Range('H5').Formula = '=''Technologies to power CCS''!D22'
    $ws4.Range('I5').Formula = '=''Technologies to power CCS''!E22'
    $ws4.Range('J5').Formula = '=''Technologies to power CCS''!F22'
    $ws4.Range('K5').Formula = '=''Technologies to power CCS''!G22'
    $ws4.Range('L5').Formula = '=''Technologies to power CCS''!H22'
    $ws4.Range('M5').Formula = '=''Technologies to power CCS''!I22'
    $ws4.Range('N5').Formula = '=''Technologies to power CCS''!J22'
    $ws4.Range('O5').Formula = '=''Technologies to power CCS''!K22'
    $ws4.Range('P5').Formula = '=''Technologies to power CCS''!L22'
    $ws4.Range('Q5').Formula = '=''Technologies to power CCS''!M22'
    $ws4.Range('R5').Formula = '=''Technologies to power CCS''!N22'
    $ws4.Range('S5').Formula = '=''Technologies to power CCS''!O22'
    $ws4.Range('T5').Formula = '=S5'
    $ws4.Range('U5').Formula = '=T5'
    $ws4.Range('V5').Formula = '=U5'
    $ws4.Range('W5').Formula = '=V5'
    $ws4.Range('X5').Formula = '=W5'
    $ws4.Range('Y5').Formula = '=X5'
    $ws4.Range('Z5').Formula = '=Y5'
    $ws4.Range('AA5').Formula = '=Z5'
    $ws4.Range('AB5').Formula = '=AA5'
    $ws4.Range('AC5').Formula = '=AB5'
    $ws4.Range('AD5').Formula = '=AC5'
    $ws4.Range('AE5').Formula = '=AD5'
    $ws4.Range('AF5').Formula = '=AE5'
    $ws4.Range('AG5').Formula = '=AF5'
    $ws4.Range('AH5').Formula = '=AG5'
    $ws4.Range('AI5').Formula = '=AH5'
    $ws4.Range('AJ5').Formula = '=AI5'
    $ws4.Range('AK5').Formula = '=AJ5'
    $ws4.Range('AL5').Formula = '=AK5'
    $ws4.Range('AM5').Formula = '=AL5'
    $ws4.Range('B6').Formula = '=C6'
    $ws4.Range('C6').Formula = '=D6'
    $ws4.Range('D6').Formula = '=E6'
    $ws4.Range('E6').Formula = '=F6'
    $ws4.Range('F6').Formula = '=''Technologies to power CCS''!B23'
    $ws4.Range('G6').Formula = '=''Technologies to power CCS''!C23'
    $ws4.Range('H6').Formula = '=''Technologies to power CCS''!D23'
    $ws4.Range('I6').Formula = '=''Technologies to power CCS''!E23'
    $ws4.Range('J6').Formula = '=''Technologies to power CCS''!F23'
    $ws4.Range('K6').Formula = '=''Technologies to power CCS''!G23'
    $ws4.Range('L6').Formula = '=''Technologies to power CCS''!H23'
    $ws4.Range('M6').Formula = '=''Technologies to power CCS''!I23'
    $ws4.Range('N6').Formula = '=''Technologies to power CCS''!J23'
    $ws4.Range('O6').Formula = '=''Technologies to power CCS''!K23'
    $ws4.Range('P6').Formula = '=''Technologies to power CCS''!L23'
    $ws4.Range('Q6').Formula = '=''Technologies to power CCS''!M23'
    $ws4.Range('R6').Formula = '=''Technologies to power CCS''!N23'
    $ws4.Range('S6').Formula = '=''Technologies to power CCS''!O23'
    $ws4.Range('T6').Formula = '=S6'
    $ws4.Range('U6').Formula = '=T6'
    $ws4.Range('V6').Formula = '=U6'
    $ws4.Range('W6').Formula = '=V6'
    $ws4.Range('X6').Formula = '=W6'
    $ws4.Range('Y6').Formula = '=X6'
    $ws4.Range('Z6').Formula = '=Y6'
    $ws4.Range('AA6').Formula = '=Z6'
    $ws4.Range('AB6').Formula = '=AA6'
    $ws4.Range('AC6').Formula = '=AB6'
    $ws4.Range('AD6').Formula = '=AC6'
    $ws4.Range('AE6').Formula = '=AD6'
    $ws4.Range('AF6').Formula = '=AE6'
    $ws4.Range('AG6').Formula = '=AF6'
    $ws4.Range('AH6').Formula = '=AG6'
    $ws4.Range('AI6').Formula = '=AH6'
    $ws4.Range('AJ6').Formula = '=AI6'
    $ws4.Range('AK6').Formula = '=AJ6'
    $ws4.Range('AL6').Formula = '=AK6'
    $ws4.Range('AM6').Formula = '=AL6'
    $ws4.Range('B7').Formula = '=C7'
    $ws4.Range('C7').Formula = '=D7'
    $ws4.Range('D7').Formula = '=E7'
    $ws4.Range('E7').Formula = '=F7'
    $ws4.Range('F7').Formula = '=''Technologies to power CCS''!B24'
    $ws4.Range('G7').Formula = '=''Technologies to power CCS''!C24'
    $ws4.Range('H7').Formula = '=''Technologies to power CCS''!D24'
    $ws4.Range('I7').Formula = '=''Technologies to power CCS''!E24'
    $ws4.Range('J7').Formula = '=''Technologies to power CCS''!F24'
    $ws4.Range('K7').Formula = '=''Technologies to power CCS''!G24'
    $ws4.Range('L7').Formula = '=''Technologies to power CCS''!H24'
    $ws4.Range('M7').Formula = '=''Technologies to power CCS''!I24'
    $ws4.Range('N7').Formula = '=''Technologies to power CCS''!J24'
    $ws4.Range('O7').Formula = '=''Technologies to power CCS''!K24'
    $ws4.Range('P7').Formula = '=''Technologies to power CCS''!L24'
    $ws4.Range('Q7').Formula = '=''Technologies to power CCS''!M24'
    $ws4.Range('R7').Formula = '=''Technologies to power CCS''!N24'
    $ws4.Range('S7').Formula = '=''Technologies to power CCS''!O24'
    $ws4.Range('T7').Formula = '=S7'
    $ws4.Range('U7').Formula = '=T7'
    $ws4.Range('V7').Formula = '=U7'
    $ws4.Range('W7').Formula = '=V7'
    $ws4.Range('X7').Formula = '=W7'
    $ws4.Range('Y7').Formula = '=X7'
    $ws4.Range('Z7').Formula = '=Y7'
    $ws4.Range('AA7').Formula = '=Z7'
    $ws4.Range('AB7').Formula = '=AA7'
    $ws4.Range('AC7').Formula = '=AB7'
    $ws4.Range('AD7').Formula = '=AC7'
    $ws4.Range('AE7').Formula = '=AD7'
    $ws4.Range('AF7').Formula = '=AE7'
    $ws4.Range('AG7').Formula = '=AF7'
    $ws4.Range('AH7').Formula = '=AG7'
    $ws4.Range('AI7').Formula = '=AH7'
    $ws4.Range('AJ7').Formula = '=AI7'
    $ws4.Range('AK7').Formula = '=AJ7'
    $ws4.Range('AL7').Formula = '=AK7'
    $ws4.Range('AM7').Formula = '=AL7'
    $ws4.Range('B8').Formula = '=C8'
    $ws4.Range('C8').Formula = '=D8'
    $ws4.Range('D8').Formula = '=E8'
    $ws4.Range('E8').Formula = '=F8'
    $ws4.Range('F8').Formula = '=''Technologies to power CCS''!B25'
    $ws4.Range('G8').Formula = '=''Technologies to power CCS''!C25'
    $ws4.Range('H8').Formula = '=''Technologies to power CCS''!D25'
    $ws4.Range('I8').Formula = '=''Technologies to power CCS''!E25'
    $ws4.Range('J8').Formula = '=''Technologies to power CCS''!F25'
    $ws4.Range('K8').Formula = '=''Technologies to power CCS''!G25'
    $ws4.Range('L8').Formula = '=''Technologies to power CCS''!H25'
    $ws4.Range('M8').Formula = '=''Technologies to power CCS''!I25'
    $ws4.Range('N8').Formula = '=''Technologies to power CCS''!J25'
    $ws4.Range('O8').Formula = '=''Technologies to power CCS''!K25'
    $ws4.Range('P8').Formula = '=''Technologies to power CCS''!L25'
    $ws4.Range('Q8').Formula = '=''Technologies to power CCS''!M25'
    $ws4.Range('R8').Formula = '=''Technologies to power CCS''!N25'
    $ws4.Range('S8').Formula = '=''Technologies to power CCS''!O25'
    $ws4.Range('T8').Formula = '=S8'
    $ws4.Range('U8').Formula = '=T8'
    $ws4.Range('V8').Formula = '=U8'
    $ws4.Range('W8').Formula = '=V8'
    $ws4.Range('X8').Formula = '=W8'
    $ws4.Range('Y8').Formula = '=X8'
    $ws4.Range('Z8').Formula = '=Y8'
    $ws4.Range('AA8').Formula = '=Z8'
    $ws4.Range('AB8').Formula = '=AA8'
    $ws4.Range('AC8').Formula = '=AB8'
    $ws4.Range('AD8').Formula = '=AC8'
    $ws4.Range('AE8').Formula = '=AD8'
    $ws4.Range('AF8').Formula = '=AE8'
    $ws4.Range('AG8').Formula = '=AF8'
    $ws4.Range('AH8').Formula = '=AG8'
    $ws4.Range('AI8').Formula = '=AH8'
    $ws4.Range('AJ8').Formula = '=AI8'
    $ws4.Range('AK8').Formula = '=AJ8'
    $ws4.Range('AL8').Formula = '=AK8'
    $ws4.Range('AM8').Formula = '=AL8'
    $ws4.Range('B9').Formula = '=C9'
    $ws4.Range('C9').Formula = '=D9'
    $ws4.Range('D9').Formula = '=E9'
    $ws4.Range('E9').Formula = '=F9'
    $ws4.Range('F9').Formula = '=''Technologies to power CCS''!B26'
    $ws4.Range('G9').Formula = '=''Technologies to power CCS''!C26'
    $ws4.Range('H9').Formula = '=''Technologies to power CCS''!D26'
    $ws4.Range('I9').Formula = '=''Technologies to power CCS''!E26'
    $ws4.Range('J9').Formula = '=''Technologies to power CCS''!F26'
    $ws4.Range('K9').Formula = '=''Technologies to power CCS''!G26'
    $ws4.Range('L9').Formula = '=''Technologies to power CCS''!H26'
    $ws4.Range('M9').Formula = '=''Technologies to power CCS''!I26'
    $ws4.Range('N9').Formula = '=''Technologies to power CCS''!J26'
    $ws4.Range('O9').Formula = '=''Technologies to power CCS''!K26'
    $ws4.Range('P9').Formula = '=''Technologies to power CCS''!L26'
    $ws4.Range('Q9').Formula = '=''Technologies to power CCS''!M26'
    $ws4.Range('R9').Formula = '=''Technologies to power CCS''!N26'
    $ws4.Range('S9').Formula = '=''Technologies to power CCS''!O26'
    $ws4.Range('T9').Formula = '=S9'
    $ws4.Range('U9').Formula = '=T9'
    $ws4.Range('V9').Formula = '=U9'
    $ws4.Range('W9').Formula = '=V9'
    $ws4.Range('X9').Formula = '=W9'
    $ws4.Range('Y9').Formula = '=X9'
    $ws4.Range('Z9').Formula = '=Y9'
    $ws4.Range('AA9').Formula = '=Z9'
    $ws4.Range('AB9').Formula = '=AA9'
    $ws4.Range('AC9').Formula = '=AB9'
    $ws4.Range('AD9').Formula = '=AC9'
    $ws4.Range('AE9').Formula = '=AD9'
    $ws4.Range('AF9').Formula = '=AE9'
    $ws4.Range('AG9').Formula = '=AF9'
    $ws4.Range('AH9').Formula = '=AG9'
    $ws4.Range('AI9').Formula = '=AH9'
    $ws4.Range('AJ9').Formula = '=AI9'
    $ws4.Range('AK9').Formula = '=AJ9'
    $ws4.Range('AL9').Formula = '=AK9'
    $ws4.Range('AM9').Formula = '=AL9'
    $ws4.Range('B10').Formula = '=C10'
    $ws4.Range('C10').Formula = '=D10'
    $ws4.Range('D10').Formula = '=E10'
    $ws4.Range('E10').Formula = '=F10'
    $ws4.Range('F10').Formula = '=''Technologies to power CCS''!B27'
    $ws4.Range('G10').Formula = '=''Technologies to power CCS''!C27'
    $ws4.Range('H10').Formula = '=''Technologies to power CCS''!D27'
    $ws4.Range('I10').Formula = '=''Technologies to power CCS''!E27'
    $ws4.Range('J10').Formula = '=''Technologies to power CCS''!F27'
    $ws4.Range('K10').Formula = '=''Technologies to power CCS''!G27'
    $ws4.Range('L10').Formula = '=''Technologies to power CCS''!H27'
    $ws4.Range('M10').Formula = '=''Technologies to power CCS''!I27'
    $ws4.Range('N10').Formula = '=''Technologies to power CCS''!J27'
    $ws4.Range('O10').Formula = '=''Technologies to power CCS''!K27'
    $ws4.Range('P10').Formula = '=''Technologies to power CCS''!L27'
    $ws4.Range('Q10').Formula = '=''Technologies to power CCS''!M27'
    $ws4.Range('R10').Formula = '=''Technologies to power CCS''!N27'
    $ws4.Range('S10').Formula = '=''Technologies to power CCS''!O27'
    $ws4.Range('T10').Formula = '=S10'
    $ws4.Range('U10').Formula = '=T10'
    $ws4.Range('V10').Formula = '=U10'
    $ws4.Range('W10').Formula = '=V10'
    $ws4.Range('X10').Formula = '=W10'
    $ws4.Range('Y10').Formula = '=X10'
    $ws4.Range('Z10').Formula = '=Y10'
    $ws4.Range('AA10').Formula = '=Z10'
    $ws4.Range('AB10').Formula = '=AA10'
    $ws4.Range('AC10').Formula = '=AB10'
    $ws4.Range('AD10').Formula = '=AC10'
    $ws4.Range('AE10').Formula = '=AD10'
    $ws4.Range('AF10').Formula = '=AE10'
    $ws4.Range('AG10').Formula = '=AF10'
    $ws4.Range('AH10').Formula = '=AG10'
    $ws4.Range('AI10').Formula = '=AH10'
    $ws4.Range('AJ10').Formula = '=AI10'
    $ws4.Range('AK10').Formula = '=AJ10'
    $ws4.Range('AL10').Formula = '=AK10'
    $ws4.Range('AM10').Formula = '=AL10'
    $ws4.Range('B11').Formula = '=C11'
    $ws4.Range('C11').Formula = '=D11'
    $ws4.Range('D11').Formula = '=E11'
    $ws4.Range('E11').Formula = '=F11'
    $ws4.Range('F11').Formula = '=''Technologies to power CCS''!B28'
    $ws4.Range('G11').Formula = '=''Technologies to power CCS''!C28'
    $ws4.Range('H11').Formula = '=''Technologies to power CCS''!D28'
    $ws4.Range('I11').Formula = '=''Technologies to power CCS''!E28'
    $ws4.Range('J11').Formula = '=''Technologies to power CCS''!F28'
    $ws4.Range('K11').Formula = '=''Technologies to power CCS''!G28'
    $ws4.Range('L11').Formula = '=''Technologies to power CCS''!H28'
    $ws4.Range('M11').Formula = '=''Technologies to power CCS''!I28'
    $ws4.Range('N11').Formula = '=''Technologies to power CCS''!J28'
    $ws4.Range('O11').Formula = '=''Technologies to power CCS''!K28'
    $ws4.Range('P11').Formula = '=''Technologies to power CCS''!L28'
    $ws4.Range('Q11').Formula = '=''Technologies to power CCS''!M28'
    $ws4.Range('R11').Formula = '=''Technologies to power CCS''!N28'
    $ws4.Range('S11').Formula = '=''Technologies to power CCS''!O28'
    $ws4.Range('T11').Formula = '=S11'
    $ws4.Range('U11').Formula = '=T11'
    $ws4.Range('V11').Formula = '=U11'
    $ws4.Range('W11').Formula = '=V11'
    $ws4.Range('X11').Formula = '=W11'
    $ws4.Range('Y11').Formula = '=X11'
    $ws4.Range('Z11').Formula = '=Y11'
    $ws4.Range('AA11').Formula = '=Z11'
    $ws4.Range('AB11').Formula = '=AA11'
    $ws4.Range('AC11').Formula = '=AB11'
    $ws4.Range('AD11').Formula = '=AC11'
    $ws4.Range('AE11').Formula = '=AD11'
    $ws4.Range('AF11').Formula = '=AE11'
    $ws4.Range('AG11').Formula = '=AF11'
    $ws4.Range('AH11').Formula = '=AG11'
    $ws4.Range('AI11').Formula = '=AH11'
    $ws4.Range('AJ11').Formula = '=AI11'
    $ws4.Range('AK11').Formula = '=AJ11'
    $ws4.Range('AL11').Formula = '=AK11'
    $ws4.Range('AM11').Formula = '=AL11'
    $ws4.Range('B12').Formula = '=C12'
    $ws4.Range('C12').Formula = '=D12'
    $ws4.Range('D12').Formula = '=E12'
    $ws4.Range('E12').Formula = '=F12'
    $ws4.Range('F12').Formula = '=''Technologies to power CCS''!B29'
    $ws4.Range('G12').Formula = '=''Technologies to power CCS''!C29'
    $ws4.Range('H12').Formula = '=''Technologies to power CCS''!D29'
    $ws4.Range('I12').Formula = '=''Technologies to power CCS''!E29'
    $ws4.Range('J12').Formula = '=''Technologies to power CCS''!F29'
    $ws4.Range('K12').Formula = '=''Technologies to power CCS''!G29'
    $ws4.Range('L12').Formula = '=''Technologies to power CCS''!H29'
    $ws4.Range('M12').Formula = '=''Technologies to power CCS''!I29'
    $ws4.Range('N12').Formula = '=''Technologies to power CCS''!J29'
    $ws4.Range('O12').Formula = '=''Technologies to power CCS''!K29'
    $ws4.Range('P12').Formula = '=''Technologies to power CCS''!L29'
    $ws4.Range('Q12').Formula = '=''Technologies to power CCS''!M29'
    $ws4.Range('R12').Formula = '=''Technologies to power CCS''!N29'
    $ws4.Range('S12').Formula = '=''Technologies to power CCS''!O29'
    $ws4.Range('T12').Formula = '=S12'
    $ws4.Range('U12').Formula = '=T12'
    $ws4.Range('V12').Formula = '=U12'
    $ws4.Range('W12').Formula = '=V12'
    $ws4.Range('X12').Formula = '=W12'
    $ws4.Range('Y12').Formula = '=X12'
    $ws4.Range('Z12').Formula = '=Y12'
    $ws4.Range('AA12').Formula = '=Z12'
    $ws4.Range('AB12').Formula = '=AA12'
    $ws4.Range('AC12').Formula = '=AB12'
    $ws4.Range('AD12').Formula = '=AC12'
    $ws4.Range('AE12').Formula = '=AD12'
    $ws4.Range('AF12').Formula = '=AE12'
    $ws4.Range('AG12').Formula = '=AF12'
    $ws4.Range('AH12').Formula = '=AG12'
    $ws4.Range('AI12').Formula = '=AH12'
    $ws4.Range('AJ12').Formula = '=AI12'
    $ws4.Range('AK12').Formula = '=AJ12'
    $ws4.Range('AL12').Formula = '=AK12'
    $ws4.Range('AM12').Formula = '=AL12'
    $ws4.Range('B13').Formula = '=C13'
    $ws4.Range('C13').Formula = '=D13'
    $ws4.Range('D13').Formula = '=E13'
    $ws4.Range('E13').Formula = '=F13'
    $ws4.Range('F13').Formula = '=''Technologies to power CCS''!B30'
    $ws4.Range('G13').Formula = '=''Technologies to power CCS''!C30'
    $ws4.Range('H13').Formula = '=''Technologies to power CCS''!D30'
    $ws4.Range('I13').Formula = '=''Technologies to power CCS''!E30'
    $ws4.Range('J13').Formula = '=''Technologies to power CCS''!F30'
    $ws4.Range('K13').Formula = '=''Technologies to power CCS''!G30'
    $ws4.Range('L13').Formula = '=''Technologies to power CCS''!H30'
    $ws4.Range('M13').Formula = '=''Technologies to power CCS''!I30'
    $ws4.Range('N13').Formula = '=''Technologies to power CCS''!J30'
    $ws4.Range('O13').Formula = '=''Technologies to power CCS''!K30'
    $ws4.Range('P13').Formula = '=''Technologies to power CCS''!L30'
    $ws4.Range('Q13').Formula = '=''Technologies to power CCS''!M30'
    $ws4.Range('R13').Formula = '=''Technologies to power CCS''!N30'
    $ws4.Range('S13').Formula = '=''Technologies to power CCS''!O30'
    $ws4.Range('T13').Formula = '=S13'
    $ws4.Range('U13').Formula = '=T13'
    $ws4.Range('V13').Formula = '=U13'
    $ws4.Range('W13').Formula = '=V13'
    $ws4.Range('X13').Formula = '=W13'
    $ws4.Range('Y13').Formula = '=X13'
    $ws4.Range('Z13').Formula = '=Y13'
    $ws4.Range('AA13').Formula = '=Z13'
    $ws4.Range('AB13').Formula = '=AA13'
    $ws4.Range('AC13').Formula = '=AB13'
    $ws4.Range('AD13').Formula = '=AC13'
    $ws4.Range('AE13').Formula = '=AD13'
    $ws4.Range('AF13').Formula = '=AE13'
    $ws4.Range('AG13').Formula = '=AF13'
    $ws4.Range('AH13').Formula = '=AG13'
    $ws4.Range('AI13').Formula = '=AH13'
    $ws4.Range('AJ13').Formula = '=AI13'
    $ws4.Range('AK13').Formula = '=AJ13'
    $ws4.Range('AL13').Formula = '=AK13'
    $ws4.Range('AM13').Formula = '=AL13'
    $ws4.Range('B14').Formula = '=C14'
    $ws4.Range('C14').Formula = '=D14'
    $ws4.Range('D14').Formula = '=E14'
    $ws4.Range('E14').Formula = '=F14'
    $ws4.Range('F14').Formula = '=''Technologies to power CCS''!B31'
    $ws4.Range('G14').Formula = '=''Technologies to power CCS''!C31'
    $ws4.Range('H14').Formula = '=''Technologies to power CCS''!D31'
    $ws4.Range('I14').Formula = '=''Technologies to power CCS''!E31'
    $ws4.Range('J14').Formula = '=''Technologies to power CCS''!F31'
    $ws4.Range('K14').Formula = '=''Technologies to power CCS''!G31'
    $ws4.Range('L14').Formula = '=''Technologies to power CCS''!H31'
    $ws4.Range('M14').Formula = '=''Technologies to power CCS''!I31'
    $ws4.Range('N14').Formula = '=''Technologies to power CCS''!J31'
    $ws4.Range('O14').Formula = '=''Technologies to power CCS''!K31'
    $ws4.Range('P14').Formula = '=''Technologies to power CCS''!L31'
    $ws4.Range('Q14').Formula = '=''Technologies to power CCS''!M31'
    $ws4.Range('R14').Formula = '=''Technologies to power CCS''!N31'
    $ws4.Range('S14').Formula = '=''Technologies to power CCS''!O31'
    $ws4.Range('T14').Formula = '=S14'
    $ws4.Range('U14').Formula = '=T14'
    $ws4.Range('V14').Formula = '=U14'
    $ws4.Range('W14').Formula = '=V14'
    $ws4.Range('X14').Formula = '=W14'
    $ws4.Range('Y14').Formula = '=X14'
    $ws4.Range('Z14').Formula = '=Y14'
    $ws4.Range('AA14').Formula = '=Z14'
    $ws4.Range('AB14').Formula = '=AA14'
    $ws4.Range('AC14').Formula = '=AB14'
    $ws4.Range('AD14').Formula = '=AC14'
    $ws4.Range('AE14').Formula = '=AD14'
    $ws4.Range('AF14').Formula = '=AE14'
    $ws4.Range('AG14').Formula = '=AF14'
    $ws4.Range('AH14').Formula = '=AG14'
    $ws4.Range('AI14').Formula = '=AH14'
    $ws4.Range('AJ14').Formula = '=AI14'
    $ws4.Range('AK14').Formula = '=AJ14'
    $ws4.Range('AL14').Formula = '=AK14'
    $ws4.Range('AM14').Formula = '=AL14'
# --- Update view selections to match target state ---
$ws3.Activate()
$ws3.Range('F37').Select()

$ws4.Activate()
$ws4.Range('H5').Select()

# --- Restore the originally active sheet/selection ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()

Write-Output "Edit complete"
